$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.853.39"
$ws.Range("E2").Value = "  -2.74%  "
# Row 3
$ws.Range("D3").Value = "1.889.98"
$ws.Range("E3").Value = "  -5.53%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.27"
$ws.Range("E5").Value = "  -1.94%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.17%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4583"
$ws.Range("E7").Value = "  -1.88%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3807"
$ws.Range("E8").Value = "  -3.70%  "
# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07719"
$ws.Range("E9").Value = "  -3.15%  "
# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9648"
$ws.Range("E10").Value = "  -4.11%  "
# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.05"
$ws.Range("E11").Value = "  -3.40%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.877.71"
$ws.Range("E12").Value = "  -7.72%  "
# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.946"
$ws.Range("E13").Value = "  -4.20%  "
# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.659"
$ws.Range("E14").Value = "  -3.74%  "
# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07047"
$ws.Range("E15").Value = "  -1.47%  "
# Row 16
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.19%  "
# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.19"
$ws.Range("E17").Value = "  -6.59%  "
# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009497"
$ws.Range("E18").Value = "  -5.13%  "
# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.60"
$ws.Range("E19").Value = "  -3.74%  "
# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.11%  "
# Row 21
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "28.813.64"
$ws.Range("E21").Value = "  -3.17%  "
# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.305"
$ws.Range("E22").Value = "  -4.34%  "
# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -3.81%  "
# Row 24
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.119.06"
$ws.Range("E24").Value = "  -7.50%  "
# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.081"
$ws.Range("E25").Value = "  -2.78%  "
# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.30"
$ws.Range("E26").Value = "  -1.11%  "
# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.05"
$ws.Range("E27").Value = "  -3.42%  "
# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.596"
$ws.Range("E28").Value = "  -6.98%  "
# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.12"
$ws.Range("E29").Value = "  -2.95%  "
# Row 30
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.809"
$ws.Range("E30").Value = "  -8.32%  "
# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09271"
$ws.Range("E31").Value = "  -1.98%  "
# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8519"
$ws.Range("E32").Value = "  -4.96%  "
# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.073"
$ws.Range("E33").Value = "  -4.47%  "
# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.229"
$ws.Range("E34").Value = "  -9.00%  "
# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.016"
$ws.Range("E35").Value = "  -5.26%  "
# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05674"
$ws.Range("E36").Value = "  -2.92%  "
# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.140"
$ws.Range("E37").Value = "  -3.66%  "
# Row 38
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  -0.08%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02032"
$ws.Range("E39").Value = "  -4.75%  "
# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.403"
$ws.Range("E40").Value = "  -6.87%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5486"
$ws.Range("E41").Value = "  -5.16%  "
# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1750"
$ws.Range("E42").Value = "  -4.12%  "
# Row 43
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000002888"
$ws.Range("E43").Value = "  -24.19%  "
# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.156"
$ws.Range("E44").Value = "  -7.56%  "
# Row 45
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.692"
$ws.Range("E45").Value = "  +1.49%  "
# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5160"
$ws.Range("E46").Value = "  -4.61%  "
# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.22"
$ws.Range("E47").Value = "  -7.70%  "
# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06805"
$ws.Range("E48").Value = "  -3.02%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.069"
$ws.Range("E49").Value = "  -5.46%  "
# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.64"
$ws.Range("E50").Value = "  -2.73%  "
# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.780"
$ws.Range("E51").Value = "  -5.24%  "
